$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.044.05'
$ws.Range("E2").Value = '  +4.55%  '
$ws.Range("D3").Value = '2.354.07'
$ws.Range("E3").Value = '  +3.15%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.39%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '519.76'
$ws.Range("E5").Value = '  +2.86%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.90'
$ws.Range("E6").Value = '  +4.84%  '
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.537'
$ws.Range("E8").Value = '  +1.65%  '
$ws.Range("D9").Value = '2.348.38'
$ws.Range("E9").Value = '  +2.15%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.104'
$ws.Range("E10").Value = '  +7.26%  '
$ws.Range("E11").Value = '  -0.74%  '
$ws.Range("E12").Value = '  +6.00%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.341'
$ws.Range("E13").Value = '  +0.06%  '
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.72'
$ws.Range("E14").Value = '  +1.10%  '
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '2.762.30'
$ws.Range("E15").Value = '  +2.69%  '
$ws.Range("D16").Value = '56.874.94'
$ws.Range("E16").Value = '  +4.13%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000135'
$ws.Range("E17").Value = '  +2.63%  '
$ws.Range("D18").Value = '2.362.62'
$ws.Range("E18").Value = '  +3.51%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.46'
$ws.Range("E19").Value = '  -0.06%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.24'
$ws.Range("E20").Value = '  +2.25%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '323.07'
$ws.Range("E21").Value = '  +5.07%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.54'
$ws.Range("E22").Value = '  -0.25%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '60.99'
$ws.Range("E24").Value = '  +1.12%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.08'
$ws.Range("E25").Value = '  +7.92%  '
$ws.Range("E26").Value = '  +7.58%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.78'
$ws.Range("E27").Value = '  +4.03%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.25'
$ws.Range("E28").Value = '  +10.38%  '
$ws.Range("E29").Value = '  -0.64%  '
$ws.Range("D30").Value = '0.0₃0746'
$ws.Range("E30").Value = '  +6.22%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.70'
$ws.Range("E31").Value = '  +3.94%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.19'
$ws.Range("E32").Value = '  +1.60%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.28'
$ws.Range("E33").Value = '  +1.64%  '
$ws.Range("E34").Value = '  +0.03%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.997'
$ws.Range("E35").Value = '  +0.30%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.25'
$ws.Range("E36").Value = '  +3.52%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.930'
$ws.Range("E37").Value = '  +2.50%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.01'
$ws.Range("E38").Value = '  +5.43%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.54'
$ws.Range("E39").Value = '  +8.54%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '37.82'
$ws.Range("E40").Value = '  +3.41%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.378'
$ws.Range("E41").Value = '  +0.72%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.60'
$ws.Range("E42").Value = '  +5.68%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '137.29'
$ws.Range("E43").Value = '  +4.72%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '278.03'
$ws.Range("E44").Value = '  +10.56%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.14'
$ws.Range("E45").Value = '  +0.85%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0936'
$ws.Range("E46").Value = '  +2.84%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0504'
$ws.Range("E47").Value = '  +1.12%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.563'
$ws.Range("E48").Value = '  +2.20%  '
$ws.Range("E49").Value = '  +5.29%  '
$ws.Range("E50").Value = '  +0.82%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.70'
$ws.Range("E51").Value = '  +11.52%  '
